$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores formatted numeric-looking text such as
# "1.00" or "21.10". Setting .Value directly would let Excel coerce
# these into real numbers and silently drop the trailing zeros, so we
# force a Text number format on each Price cell we touch before writing
# the new value, keeping it stored as text exactly as scraped.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D13", "D14", "D15", "D16", "D18", "D19", "D21", "D23", "D24", "D27", "D29", "D32", "D33", "D34", "D37", "D39", "D40", "D42", "D44", "D45", "D47")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '57.922.75'
$ws.Range("E2").Value = '  +2.21%  '

$ws.Range("D3").Value = '3.065.27'
$ws.Range("E3").Value = '  +2.74%  '

$ws.Range("D5").Value = '526.75'
$ws.Range("E5").Value = '  +5.86%  '

$ws.Range("D6").Value = '143.32'
$ws.Range("E6").Value = '  +5.63%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +5.47%  '

$ws.Range("D9").Value = '7.64'
$ws.Range("E9").Value = '  +5.86%  '

$ws.Range("E10").Value = '  +7.23%  '

$ws.Range("E11").Value = '  +5.61%  '

$ws.Range("E12").Value = '  +2.03%  '

$ws.Range("D13").Value = '3.592.59'
$ws.Range("E13").Value = '  +2.85%  '

$ws.Range("D14").Value = '27.36'
$ws.Range("E14").Value = '  +8.22%  '

$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +16.40%  '

$ws.Range("D16").Value = '57.931.65'
$ws.Range("E16").Value = '  +2.35%  '

$ws.Range("E17").Value = '  +7.05%  '

$ws.Range("D18").Value = '3.069.42'
$ws.Range("E18").Value = '  +2.98%  '

$ws.Range("D19").Value = '13.26'
$ws.Range("E19").Value = '  +7.14%  '

$ws.Range("E20").Value = '  +5.11%  '

$ws.Range("D21").Value = '341.68'
$ws.Range("E21").Value = '  +4.69%  '

$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Value = '0.503'
$ws.Range("E23").Value = '  +6.93%  '

$ws.Range("D24").Value = '64.97'
$ws.Range("E24").Value = '  +5.30%  '

$ws.Range("E25").Value = '  +8.72%  '

$ws.Range("E26").Value = '  +5.46%  '

$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.33%  '

$ws.Range("E28").Value = '  +7.27%  '

$ws.Range("D29").Value = '7.42'
$ws.Range("E29").Value = '  +9.50%  '

$ws.Range("E30").Value = '  +6.54%  '

$ws.Range("E31").Value = '  +6.65%  '

$ws.Range("D32").Value = '21.10'
$ws.Range("E32").Value = '  +3.81%  '

$ws.Range("D33").Value = '4.79'
$ws.Range("E33").Value = '  +7.34%  '

$ws.Range("D34").Value = '156.89'
$ws.Range("E34").Value = '  +2.79%  '

$ws.Range("E35").Value = '  +7.07%  '

$ws.Range("E36").Value = '  +3.79%  '

$ws.Range("D37").Value = '26.41'
$ws.Range("E37").Value = '  +13.29%  '

$ws.Range("E38").Value = '  +4.98%  '

$ws.Range("D39").Value = '3.101.66'
$ws.Range("E39").Value = '  +2.89%  '

$ws.Range("D40").Value = '37.90'
$ws.Range("E40").Value = '  +3.47%  '

$ws.Range("E41").Value = '  +10.24%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.48'
$ws.Range("E42").Value = '  +5.35%  '

$ws.Range("E43").Value = '  +0.22%  '

$ws.Range("B44").Value = 'Mantle'
$ws.Range("C44").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D44").Value = '0.667'
$ws.Range("E44").Value = '  +4.03%  '

$ws.Range("D45").Value = '2.340.91'
$ws.Range("E45").Value = '  +5.09%  '

$ws.Range("E46").Value = '  +3.32%  '

$ws.Range("D47").Value = '2.00'
$ws.Range("E47").Value = '  +2.53%  '

$ws.Range("E48").Value = '  +5.29%  '

$ws.Range("E49").Value = '  +3.81%  '

$ws.Range("E50").Value = '  +6.00%  '

$ws.Range("E51").Value = '  +6.17%  '
